$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Entry order chosen to reproduce the original shared-string table order
$ws.Range("A3").Value = "Bootbot"

$ws.Range("B2").Value = "Price"
$ws.Range("C2").Value = "Response time"
$ws.Range("D2").Value = "Languages used"
$ws.Range("E2").Value = "Extendibility"
$ws.Range("F2").Value = "Range of platforms"

$ws.Range("B3").Value = "Free"

$ws.Range("A4").Value = "Bootbot scores /5"

$ws.Range("D3").Value = "Javascript"

$ws.Range("A5").Value = "Dialog flow"

$ws.Range("A6").Value = "Pandora bot"

$ws.Range("B4").Value = 5

# Column widths (values chosen so the engine's pixel-quantized ColumnWidth
# storage reproduces the target XML "width" attribute as closely as possible)
$ws.Range("A1").ColumnWidth = 17.333333333333332
$ws.Range("C1").ColumnWidth = 13.166666666666666
$ws.Range("D1").ColumnWidth = 14.166666666666666
$ws.Range("E1").ColumnWidth = 11.5
$ws.Range("F1").ColumnWidth = 17.666666666666668

# Selection
$ws.Range("A6").Select()
